# This edit adds one new weekly price observation for "Cebollín" at
# "Terminal Hortofrutícola Agro Chillán" (Ñuble). The new record is inserted
# as a new row 138, pushing the previously existing rows 138-188 down to
# rows 139-189 (dimension grows from A1:R188 to A1:R189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 138 (shifts rows 138:188 down to 139:189,
# and the new row inherits formatting - e.g. the date number format - from the
# row above it, matching the style="2" seen on column D in the target sheet).
$ws.Rows(138).Insert()

# Populate the newly inserted row 138 with the new observation's data.
$ws.Range("A138").Value = 7
$ws.Range("B138").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C138").Value = "Ñuble"
$ws.Range("D138").Value = 45141
$ws.Range("E138").Value = 16
$ws.Range("F138").Value = 100112037
$ws.Range("G138").Value = "Cebollín"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 100
$ws.Range("K138").Value = 6500
$ws.Range("L138").Value = 6500
$ws.Range("M138").Value = 6500
$ws.Range("N138").Value = "$/paquete 36 unidades"
$ws.Range("O138").Value = "Provincia de Diguillín"
$ws.Range("P138").Value = 181
$ws.Range("Q138").Value = 36
$ws.Range("R138").Value = "Hortaliza"
